# Tee.xlsx update: add "Average Elevation (m)" column, fix a couple of
# strings, and append the "Mate" tea row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Append the new "Mate" row (row 19) using the CURRENT (pre-insert)
#    column layout: H = Type, I = Color. Copy the formatting from the
#    row above first so number formats / fonts match the rest of the
#    table.
# ---------------------------------------------------------------------
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)

$ws.Range("A19").Value = "Mate"
$ws.Range("B19").Value = 0.01
$ws.Range("C19").Value = 30.0
$ws.Range("D19").Formula = "=B19/(C19+B19)"
$ws.Range("E19").Value = "Argentina, Brazil"
$ws.Range("F19").Value = 400.0
$ws.Range("G19").Value = 800.0
$ws.Range("H19").Value = "green"
$ws.Range("I19").Value = "green"

# ---------------------------------------------------------------------
# 2) Insert a new column at H. This pushes the old H ("Type") values to
#    I and the old I ("Color") values to J, for every row including the
#    new row 19 - exactly matching the desired final layout.
# ---------------------------------------------------------------------
$ws.Columns.Item(8).Insert()

# ---------------------------------------------------------------------
# 3) Fill the new H column with the "Average Elevation (m)" header and
#    the average-of-min/max formula.
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "Average Elevation (m)"
$ws.Range("H2:H19").Formula = "=(F2+G2)/2"

# Widen column H (closest value reachable through the pixel-snapped
# ColumnWidth property to the authored width of 19.38).
$ws.Columns.Item(8).ColumnWidth = 18.5

# ---------------------------------------------------------------------
# 4) Small text fixes elsewhere in the sheet.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Earl Grey"        # "Earl grey" -> "Earl Grey"
$ws.Range("E4").Value = "China, Japan"     # Sencha location correction

Write-Host "done"
